$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Reword the existing "Create billingSlice.js" action item
$ws.Range("A7").Value = "Action: Create billingSlice.js"

# 2. Open up the row skeleton for the two new "Action:" bullets (and their
#    surrounding blank spacer rows), shifting the numbered checklist items down.
$ws.Rows("8:9").Insert()
$ws.Rows("11:13").Insert()
$ws.Rows("15:16").Insert()

# 3. New action item under "2. Create applicable Stripe objects..." (row 10)
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Action: map existing users in Admin service to Stripe Customer objects"

# 4. New action item under "3. Teachers Can Create Customer" (row 14)
$ws.Range("A10").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Action: Generate a new customer when a user sets up billing for a student"

# 5. Row-height tweaks that came along with the edits made further down the sheet
$ws.Rows("13:29").RowHeight = 17

# 6. Final selection left on the sheet
$ws.Range("L18").Select()
